$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to text so numeric-looking strings
# (e.g. "1.00", "578.64") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.914.57"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "3.415.05"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "578.64"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").Value = "144.35"
$ws.Range("E6").Value = "  +2.38%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "7.59"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "4.002.37"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "28.44"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "3.413.14"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "61.952.03"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").Value = "6.16"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").Value = "13.97"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D21").Value = "390.43"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").Value = "74.85"
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.556.29"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "0.0000114"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("D28").Value = "7.47"
$ws.Range("E28").Value = "  +2.89%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.40"
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("D34").Value = "23.55"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").Value = "5.26"
$ws.Range("E35").Value = "  +5.74%  "
$ws.Range("D36").Value = "6.98"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").Value = "167.50"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").Value = "1.53"
$ws.Range("E38").Value = "  +4.76%  "
$ws.Range("D39").Value = "3.447.76"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").Value = "28.58"
$ws.Range("E40").Value = "  +8.57%  "
$ws.Range("E41").Value = "  -1.49%  "
$ws.Range("D42").Value = "0.784"
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").Value = "1.17"
$ws.Range("E45").Value = "  +4.39%  "
$ws.Range("D46").Value = "2.504.07"
$ws.Range("E46").Value = "  +2.27%  "
$ws.Range("D47").Value = "22.83"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "6.63"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("E51").Value = "  -2.43%  "

# Restore column D style/format to its original (unstyled) state.
$ws.Range("D2:D51").Style = "Normal"
